$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arbeitsmatrix")

# ------------------------------------------------------------------
# 1. Insert three new rows starting at row 155. This pushes the old
#    row 155 (blank placeholder) down to row 158, and the old summary
#    block (rows 156-159) down to rows 159-162.
# ------------------------------------------------------------------
$ws.Rows.Item(155).Insert()
$ws.Rows.Item(155).Insert()
$ws.Rows.Item(155).Insert()

# ------------------------------------------------------------------
# 2. Copy cell formatting from existing similar rows so the new rows
#    155/156 pick up the same styles Excel would use for a "last
#    entry of a multi-entry day" row (right aligned prefix column,
#    tinted time columns, etc. - same pattern as rows 139-142).
# ------------------------------------------------------------------
$ws.Range("A142:K142").Copy()
$ws.Range("A155:K156").PasteSpecial(-4122)
$ws.Range("J139").Copy()
$ws.Range("J156").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 3. Fill in the two new data rows.
# ------------------------------------------------------------------
# Row 155: "Tablet und Mobile Rezept ansicht"
$ws.Cells.Item(155, 1).Value = 22
$ws.Cells.Item(155, 2).Value = "Interface Design"
$ws.Cells.Item(155, 3).Value = "MockUps"
$ws.Cells.Item(155, 4).Value = "[FEATURE]"
$ws.Cells.Item(155, 5).Value = "Tablet und Mobile Rezept ansicht"
$ws.Cells.Item(155, 6).Value = 44501
$ws.Cells.Item(155, 7).Value = 44481
$ws.Cells.Item(155, 9).Formula = "=ROUNDUP(((SUM(K155-J155)*24*60/60)/0.25),0)*0.25"
$ws.Cells.Item(155, 10).Value = 0.45833333333333331
$ws.Cells.Item(155, 11).Value = 0.60416666666666663

# Row 156: "Tablet und Mobile Profile"
$ws.Cells.Item(156, 1).Value = 22
$ws.Cells.Item(156, 2).Value = "Interface Design"
$ws.Cells.Item(156, 3).Value = "MockUps"
$ws.Cells.Item(156, 4).Value = "[FEATURE]"
$ws.Cells.Item(156, 5).Value = "Tablet und Mobile Profile"
$ws.Cells.Item(156, 6).Value = 44501
$ws.Cells.Item(156, 7).Value = 44481
$ws.Cells.Item(156, 9).Formula = "=ROUNDUP(((SUM(K156-J156)*24*60/60)/0.25),0)*0.25"
$ws.Cells.Item(156, 10).Formula = "=K155"
$ws.Cells.Item(156, 11).Value = 0.69791666666666663

# ------------------------------------------------------------------
# 4. The two new rows already carry an explicit "[FEATURE]" prefix
#    (like the other "last entry of the day" rows 114/142), so they
#    should not be part of the Prefix dropdown-list validation -
#    remove it there, same as Excel would leave it off a cell whose
#    format no longer matches the validated range.
# ------------------------------------------------------------------
$ws.Range("D155:D156").Validation.Delete()

# ------------------------------------------------------------------
# 5. Leave the selection on the newly added cell, like the author
#    did after finishing data entry.
# ------------------------------------------------------------------
$ws.Range("G155").Select()
